# arreglo de variable contador
# Move the "VOLTAJE 1" / "VOLTAJE 2" counter-variable labels from column J
# (rows 9-10) down to column E (rows 14-15), keeping their original
# formatting, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move J9 ("VOLTAJE 1") -> E14, copying value + formatting, then clear source.
$ws.Range("J9").Copy($ws.Range("E14"))
$ws.Range("J9").Clear()

# Move J10 ("VOLTAJE 2") -> E15, copying value + formatting, then clear source.
$ws.Range("J10").Copy($ws.Range("E15"))
$ws.Range("J10").Clear()

# Update the selection to match the new editing position.
$ws.Range("H19").Select()
